# Aula 07 - Algoritmos e Complexidade - Grafos
# Rename the "Contextualização" title text to "Grafos" on the section-cover
# slide (slide 2) and on the slide's own title (slide 3).

$p = $ppt.ActivePresentation

# --- Slide 2: section-cover title shape ("Título 2", shape id 3) ---
# Paragraph runs (split by <a:br/>): 1 "Aula 07", 2 tab, 3 "Contextualização"
$s2  = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(3)
$tr2 = $shp2.TextFrame.TextRange
$tr2.Runs(3).Text = "Grafos"

# --- Slide 3: slide title shape ("Title 1", shape id 6) ---
$s3  = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(1)
$tr3 = $shp3.TextFrame.TextRange
$tr3.Runs(1).Text = "Grafos"
